$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# Insert a new row at position 2, shifting existing rows down
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with id_credito_testeado / 168
$ws.Range("A2").Value = "id_credito_testeado"
$ws.Range("B2").Value = 168

# Match formatting of the surrounding rows (A column label style, B column plain)
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the 'entidades' summary blob (now row 9) with re-ordered keys
$ws.Range("B9").Value = '{''ner-general'': {''entidades_correctas'': 0, ''entidades_incorrectas'': 0, ''porcentaje_entidades_correctas'': 1, ''paginas_correctas'': 0, ''paginas_incorrectas'': 0, ''porcentaje_paginas_correctas'': 1}, ''ner-escritura'': {''entidades_correctas'': 18, ''entidades_incorrectas'': 20, ''porcentaje_entidades_correctas'': 0.47368421052631576, ''paginas_correctas'': 18, ''paginas_incorrectas'': 20, ''porcentaje_paginas_correctas'': 0.8157894736842105, ''CIUDAD_INMUEBLE'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': False, ''pagina_real'': 53, ''texto_es_correcto'': False, ''valor_encontrado'': '' '', ''valor_real'': ''LEON\xa0''}, ''CLAUSULA_ACTUALIZACION_CREDITO'': {''pagina_encontrado'': 30, ''pagina_es_correcta'': True, ''pagina_real'': 30, ''texto_es_correcto'': True, ''valor_encontrado'': ''DECIMA PRIMERA'', ''valor_real'': ''DECIMA PRIMERA''}, ''CLAUSULA_CONSTITUCION_HIPOTECA'': {''pagina_encontrado'': 12, ''pagina_es_correcta'': True, ''pagina_real'': 12, ''texto_es_correcto'': True, ''valor_encontrado'': ''SEGUNDA'', ''valor_real'': ''SEGUNDA''}, ''CLAUSULA_JURISDICCION_COMPETENCIA'': {''pagina_encontrado'': 13, ''pagina_es_correcta'': True, ''pagina_real'': 13, ''texto_es_correcto'': True, ''valor_encontrado'': ''CUARTA'', ''valor_real'': ''CUARTA''}, ''CLAUSULA_PRORROGA'': {''pagina_encontrado'': 31, ''pagina_es_correcta'': True, ''pagina_real'': 31, ''texto_es_correcto'': True, ''valor_encontrado'': ''DECIMA QUINTA'', ''valor_real'': ''DÉCIMA QUINTA''}, ''CLAUSULA_VENCIMIENTO_ANTICIPADO'': {''pagina_encontrado'': 34, ''pagina_es_correcta'': True, ''pagina_real'': 34, ''texto_es_correcto'': True, ''valor_encontrado'': ''VIGESIMA PRIMERA'', ''valor_real'': ''VIGÉSIMA PRIMERA''}, ''CLAUSULA_INT_ORD'': {''pagina_encontrado'': 29, ''pagina_es_correcta'': True, ''pagina_real'': 30, ''texto_es_correcto'': True, ''valor_encontrado'': ''DECIMA'', ''valor_real'': ''DECIMA''}, ''CLAUSULA_PLAZO'': {''pagina_encontrado'': 12, ''pagina_es_correcta'': False, ''pagina_real'': 28, ''texto_es_correcto'': False, ''valor_encontrado'': ''CUARTA'', ''valor_real'': ''OCTAVA''}, ''CREDITO_PESOS'': {''pagina_encontrado'': 36, ''pagina_es_correcta'': True, ''pagina_real'': 36, ''texto_es_correcto'': False, ''valor_encontrado'': ''$340,999.81 CIENTOS. CUARENTA MIL NOVECIENTOS NOVENTA Y NUEVE PESOS 81/100 M.N.'', ''valor_real'': ''$340,999.81 MONTO (EN PESOS) DEL CREDITO OTORGADO TREŚCIENTOS. CUARENTA MIL NOVECIENTOS NOVENTA Y NUEVE PESOS 81/100 M.N.''}, ''FOLIO_REAL'': {''pagina_encontrado'': 6, ''pagina_es_correcta'': True, ''pagina_real'': 6, ''texto_es_correcto'': True, ''valor_encontrado'': ''FOLIO REAL R20*309941'', ''valor_real'': ''FOLIO REAL R20*309941''}, ''DESCRIPCION_INMUEBLE'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': True, ''pagina_real'': 3, ''texto_es_correcto'': True, ''valor_encontrado'': ''LOTE DE TERRENO NUMERO 2 DOS, DE LA MANZANA NUMERO 118 CIENTO DIECIOCHO'', ''valor_real'': ''LOTE DE TERRENO NUMERO 2 DOS, DE LA MANZANA NUMERO 118 CIENTO DIECIOCHO''}, ''DOMICILIO_GARANTIA'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': True, ''pagina_real'': 3, ''texto_es_correcto'': False, ''valor_encontrado'': ''CASA HABITACION MARCADA CON EL NUMERO 203 DOSCIENTOS TRES, UBICADA EN CALLE CIRCUITO CAMPESTRE SAUCE, CONSTRUIDA SOBRE EL LOTE DE TERRENO NUMERO 2 DOS, DE LA MANZANA NUMERO 118 CIENTO DIECIOCHO, DEL FRACCIONAMIENTO URBIVILLA DEL ROBLE DE ESTA CIUDAD DE LEON GUANAJUATO'', ''valor_real'': ''CASA HABITACION MARCADA CON EL NUMERO 203 DOSCIENTOS TRES, UBICADA EN CALLE CIRCUITO CAMPESTRE SAUCE, CONSTRUIDA SOBRE EL LOTE DE TERRENO NUMERO 2 DOS, DE LA MANZANA NUMERO 118 CIENTO DIECIOCHO, DEL FRACCIONAMIENTO URBIVILLA DEL ROBLE DE ESTA CIUDAD DE LEON GUANAJUATO.''}, ''ESTADO_INMUEBLE'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': False, ''pagina_real'': 53, ''texto_es_correcto'': True, ''valor_encontrado'': ''GUANAJUATO'', ''valor_real'': ''GUANAJUATO''}, ''FECHA_ESCRITURA'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': True, ''pagina_real'': 3, ''texto_es_correcto'': True, ''valor_encontrado'': ''07 SIETE DIAS DEL MES DE ENERO DEL AÑO 2013 DOS MIL TRECE'', ''valor_real'': ''07 SIETE DIAS DEL MES DE ENERO DEL AÑO 2013 DOS MIL TRECE''}, ''JURISDICCION_NOTARIO'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': True, ''pagina_real'': 3, ''texto_es_correcto'': True, ''valor_encontrado'': ''LEON, ESTADO DE GUANAJUATO'', ''valor_real'': ''LEON, ESTADO DE GUANAJUATO''}, ''JURISDICCION_SOMETIMIENTO_EXPRESO'': {''pagina_encontrado'': 13, ''pagina_es_correcta'': True, ''pagina_real'': 13, ''texto_es_correcto'': True, ''valor_encontrado'': ''DISTRITO FEDERAL O A LA DE LOS TRIBUNALES COMPETENTES EN EL LUGAR EN DONDE SE UBIQUE EL INMUEBLE OBJETO DE ESTA ESCRITURA'', ''valor_real'': ''DISTRITO FEDERAL O A LA DE LOS TRIBUNALES COMPETENTES EN EL LUGAR EN DONDE SE UBIQUE EL INMUEBLE OBJETO DE ESTA ESCRITURA''}, ''COLINDANCIAS_1'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': True, ''pagina_real'': 3, ''texto_es_correcto'': False, ''valor_encontrado'': ''AL NORESTE: 5.925 MTS, CINCO METROS NOVECIENTOS VEINTICINCO MILIMETROS, CON CALLE CIRCUITO CAMPESTRE SAUCE.-'', ''valor_real'': ''AL NORESTE: 5.925 MTS, CINCO METROS NOVECIENTOS VEINTICINCO MILIMETROS, CON CALLE CIRCUITO CAMPESTRE SAUCE''}, ''COLINDANCIAS_2'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': True, ''pagina_real'': 3, ''texto_es_correcto'': False, ''valor_encontrado'': ''AL SUROESTE: 5.925 MTS, CINCO METROS NOVECIENTOS VEINTICINCO MILIMETROS Y LINDA CON LA CASA 104 CIENTO CUATRO.-'', ''valor_real'': ''AL SUROESTE: 5.925 MTS, CINCO METROS NOVECIENTOS VEINTICINCO MILIMETROS Y LINDA CON LA CASA 104 CIENTO CUATRO''}, ''COLINDANCIAS_3'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': True, ''pagina_real'': 3, ''texto_es_correcto'': False, ''valor_encontrado'': ''AL SURESTE: 15.00 MTS, QUINCE METROS Y LINDA CON CASA 201-A DOSCIENTOS UNO LETRA A -'', ''valor_real'': ''AL SURESTE: 15.00 MTS, QUINCE METROS Y LINDA CON CASA 201-A DOSCIENTOS UNO LETRA A\xa0''}, ''COLINDANCIAS_4'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': True, ''pagina_real'': 3, ''texto_es_correcto'': False, ''valor_encontrado'': ''AL NOROESTE: 15.00 MTS., QUINCE METROS Y LINDA CON LA CASA 203-A DOSCIENTOS TRES LETRA A. - CON UN PRO-INDIVISO DE 50% CINCUENTA POR CIENTO. - DICHO INMUEBLE SE ENCUENTRA REGISTRADO EN EL REGISTRO DE PROPIEDAD Y DEL COMERCIO DE ESTA CIUDAD DE LEON, GUANAJUATO; BAJO EL FOLIO REAL NUMERO R20*429215.- SEGUNDO.- FIDEICOMISO. MEDIANTE ESCRITURA PUBLICA NUMERO.- 27,797 VEINTISIETE MIL SETECIENTOS NOVENTA Y SIETE DE FECHA PRIMERO DE DICIEMBRE DEL AÑO 2009 DOS MIL NUEVE, PASADA ANTE LA FE DEL LICENCIADO JORGE HUMBERTO CARPIO MENDOZA, NOTARIO PUBLICO NUMERO 95 NOVENTA Y CINCO, EN LEGAL EJERCICIO EN'', ''valor_real'': ''AL NOROESTE: 15.00 MTS., QUINCE METROS Y LINDA CON LA CASA 203-A DOSCIENTOS TRES LETRA A''}, ''MENSUALIDADES_CUBRIR'': {''pagina_encontrado'': 12, ''pagina_es_correcta'': False, ''pagina_real'': 28, ''texto_es_correcto'': False, ''valor_encontrado'': ''3 (TRES)'', ''valor_real'': ''360 (TRESCIENTAS SESENTA''}, ''NOTARIO'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': True, ''pagina_real'': 3, ''texto_es_correcto'': False, ''valor_encontrado'': '' '', ''valor_real'': ''LUIS ANGEL ALFONSO CHICO GONZALEZ''}, ''NUM_ CLAUSULA_DESCUENTO_SALARIO_REA'': {''pagina_encontrado'': 28, ''pagina_es_correcta'': True, ''pagina_real'': 28, ''texto_es_correcto'': True, ''valor_encontrado'': ''NOVENA'', ''valor_real'': ''NOVENA''}, ''NUM_ CLAUSULA_DESCUENTO_SALARIO_ROA'': {''pagina_encontrado'': 28, ''pagina_es_correcta'': True, ''pagina_real'': 28, ''texto_es_correcto'': True, ''valor_encontrado'': ''NOVENA'', ''valor_real'': ''NOVENA''}, ''NUM_CLAUSULA_INTERES_MORATORIO'': {''pagina_encontrado'': 24, ''pagina_es_correcta'': False, ''pagina_real'': 30, ''texto_es_correcto'': False, ''valor_encontrado'': ''SEGUNDA'', ''valor_real'': ''DÉCIMA SEGUNDA''}, ''NUM_CLAUSULA_OTORGAMIENTO_CREDITO'': {''pagina_encontrado'': 11, ''pagina_es_correcta'': True, ''pagina_real'': 11, ''texto_es_correcto'': True, ''valor_encontrado'': ''PRIMERA'', ''valor_real'': ''PRIMERA''}, ''NUM_NOTARIO'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': True, ''pagina_real'': 3, ''texto_es_correcto'': False, ''valor_encontrado'': ''5 CINCO'', ''valor_real'': ''\xa05 CINCO''}, ''PLAZO_ANIOS'': {''pagina_encontrado'': 12, ''pagina_es_correcta'': False, ''pagina_real'': 28, ''texto_es_correcto'': False, ''valor_encontrado'': '' '', ''valor_real'': ''30 (TREINTA) AÑOS''}, ''REA_VSM'': {''pagina_encontrado'': 36, ''pagina_es_correcta'': True, ''pagina_real'': 36, ''texto_es_correcto'': False, ''valor_encontrado'': ''1.669'', ''valor_real'': ''1.996\xa0''}, ''PORCENTAJE_ROA'': {''pagina_encontrado'': 28, ''pagina_es_correcta'': True, ''pagina_real'': 28, ''texto_es_correcto'': True, ''valor_encontrado'': ''20% (VEINTE POR CIENTO)'', ''valor_real'': ''20% (VEINTE POR CIENTO)''}, ''ROA_VSM'': {''pagina_encontrado'': 36, ''pagina_es_correcta'': True, ''pagina_real'': 36, ''texto_es_correcto'': False, ''valor_encontrado'': ''0.04'', ''valor_real'': ''1.669\xa0''}, ''SUPERFICIE_CONSTRUCCION'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': True, ''pagina_real'': 3, ''texto_es_correcto'': True, ''valor_encontrado'': ''88.875 M2. OCHENTA Y OCHO METROS OCHOCIENTOS SETENTA Y CINCO CENTIMETROS CUADRADOS'', ''valor_real'': ''88.875 M2. OCHENTA Y OCHO METROS OCHOCIENTOS SETENTA Y CINCO CENTIMETROS CUADRADOS''}, ''SUPERFICIE_TERRENO'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': True, ''pagina_real'': 3, ''texto_es_correcto'': False, ''valor_encontrado'': ''50% CINCUENTA POR CIENTO'', ''valor_real'': ''88.875 M2. OCHENTA Y OCHO METROS OCHOCIENTOS SETENTA Y CINCO CENTIMETROS CUADRADOS''}, ''TASA_INTERES_MORATORIO'': {''pagina_encontrado'': 36, ''pagina_es_correcta'': True, ''pagina_real'': 36, ''texto_es_correcto'': False, ''valor_encontrado'': ''13.7%'', ''valor_real'': ''9.5 + 4.2''}, ''TASA_INTERES_ORDINARIO'': {''pagina_encontrado'': 36, ''pagina_es_correcta'': True, ''pagina_real'': 36, ''texto_es_correcto'': False, ''valor_encontrado'': ''9.5%'', ''valor_real'': ''9.5%\xa0''}, ''TEXTO_ACTUALIZACION'': {''pagina_encontrado'': 30, ''pagina_es_correcta'': True, ''pagina_real'': 30, ''texto_es_correcto'': True, ''valor_encontrado'': ''EL TRABAJADOR ACEPTA EN ESTE ACTO QUE EL SALDO INSOLUTO DEL CREDITO SE INCREMENTARA EN LA MISMA PROPORCION EN QUE AUMENTE EL SALARIO MINIMO GENERAL DIARIO QUE RIJA EN EL DISTRITO FEDERAL, A PARTIR DE LA FECHA EN QUE ENTRE EN VIGOR EN EL DISTRITO FEDERAL EL AUMENTO DE QUE SE TRATE DE DICHO SALARIO, Y, POR CONSIGUIENTE, EL TRABAJADOR SE OBLIGA A PAGAR AL INFONAVIT EL SALDO INSOLUTO DEL CREDITO QUE RESULTARE DEL AJUSTE QUE SE HICIERE CONFORME A LO ANTES ESTIPULADO'', ''valor_real'': ''EL TRABAJADOR ACEPTA EN ESTE ACTO QUE EL SALDO INSOLUTO DEL CREDITO SE INCREMENTARA EN LA MISMA PROPORCION EN QUE AUMENTE EL SALARIO MINIMO GENERAL DIARIO QUE RIJA EN EL DISTRITO FEDERAL, A PARTIR DE LA FECHA EN QUE ENTRE EN VIGOR EN EL DISTRITO FEDERAL EL AUMENTO DE QUE SE TRATE DE DICHO SALARIO, Y, POR CONSIGUIENTE, EL TRABAJADOR SE OBLIGA A PAGAR AL INFONAVIT EL SALDO INSOLUTO DEL CREDITO QUE RESULTARE DEL AJUSTE QUE SE HICIERE CONFORME A LO ANTES ESTIPULADO''}, ''TEXTO_CLAUSULA_HIPOTECA'': {''pagina_encontrado'': 12, ''pagina_es_correcta'': True, ''pagina_real'': 12, ''texto_es_correcto'': False, ''valor_encontrado'': ''EL TRABAJADOR LA SENORA ANA KARINA MARTINEZ GONZALEZ, PARA GARANTIZAR EL CUMPLIMIENTO DE TODAS Y CADA UNA DE LAS OBLIGACIONES QUE EL TRABAJADOR CONTRAE POR VIRTUD DEL PRESENTE CONTRATO, CONSTITUYE HIPOTECA A FAVOR DEL INFONAVIT SOBRE EL INMUEBLE RELACIONADO EN EL ANTECEDENTE PRIMERO DE ESTA ESCRITURA'', ''valor_real'': ''EL TRABAJADOR LA SEÑORA ANA KARINA MARTINEZ GONZALEZ, PARA GARANTIZAR EL CUMPLIMIENTO DE TODAS Y CADA UNA DE LAS OBLIGACIONES QUE EL TRABAJADOR CONTRAE POR VIRTUD DEL PRESENTE CONTRATO, CONSTITUYE HIPOTECA A FAVOR DEL INFONAVIT SOBRE EL INMUEBLE RELACIONADO EN EL ANTECEDENTE PRIMERO DE ESTA ESCRITURA, CON LA SUPERFICIE, MEDIDAS Y LINDEROS AHI DETERMINADOS, LOS CUALES SE TIENEN AQUI POR REPRODUCIDOS COMO SI A LA LETRA SE INSERTASEN.''}, ''TIPO_INMUEBLE'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': False, ''pagina_real'': 53, ''texto_es_correcto'': False, ''valor_encontrado'': ''CASA HABITACION'', ''valor_real'': ''CASA EN CONDOMINIO''}}, ''ner-certificado'': {''entidades_correctas'': 8, ''entidades_incorrectas'': 4, ''porcentaje_entidades_correctas'': 0.6666666666666666, ''paginas_correctas'': 8, ''paginas_incorrectas'': 4, ''porcentaje_paginas_correctas'': 1.0, ''CONDICIONES _CREDITO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': True, ''valor_encontrado'': ''VSM'', ''valor_real'': ''VSM''}, ''CURP_ACREDITADO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': True, ''valor_encontrado'': ''MAGA841230MGTRNN 01'', ''valor_real'': ''MAGA841230MGTRNN 01''}, ''DELEGACION'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': True, ''valor_encontrado'': ''GUANAJUATO'', ''valor_real'': ''GUANAJUATO''}, ''FECHA_EMISION'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': False, ''valor_encontrado'': '' '', ''valor_real'': ''29 DE JULIO\xa0DE 2017''}, ''MESES_ADEUDADOS_CERTIFICADO_ADEUDOS'': {''pagina_encontrado'': 3, ''pagina_es_correcta'': True, ''pagina_real'': 3, ''texto_es_correcto'': True, ''valor_encontrado'': ''38'', ''valor_real'': ''38''}, ''MONTO_VSM_INTERES_MORATORIO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': False, ''valor_encontrado'': ''62.9330'', ''valor_real'': ''629330''}, ''NOMBRE_ACREDITADO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': True, ''valor_encontrado'': ''MARTINEZ GONZALEZ ANA KARINA'', ''valor_real'': ''MARTINEZ GONZALEZ ANA KARINA''}, ''NUMERO_CREDITO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': True, ''valor_encontrado'': ''1113002217'', ''valor_real'': ''1113002217''}, ''NUMERO_HOJAS_CA'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': True, ''valor_encontrado'': ''3'', ''valor_real'': ''3''}, ''PRORROGAS_APLICADAS'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': False, ''valor_encontrado'': ''31/03/2014'', ''valor_real'': ''41729''}, ''RFC_ACREDITADO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': True, ''valor_encontrado'': ''MAGA841230641'', ''valor_real'': ''MAGA841230641''}, ''SUERTE_PRINCIPAL'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': False, ''valor_encontrado'': ''187.7420'', ''valor_real'': ''1877420''}}, ''avaluo_adjudicacion'': {''entidades_correctas'': 0, ''entidades_incorrectas'': 0, ''porcentaje_entidades_correctas'': 1, ''paginas_correctas'': 0, ''paginas_incorrectas'': 0, ''porcentaje_paginas_correctas'': 1}}'

Write-Host "done"
